# Auto-generated edit script applying numeric corrections to the Leve profit
# tracking tables across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 116.58
$ws.Range("I15").Value = 116.58
$ws.Range("K15").Value = 349.74
$ws.Range("M15").Value = -180.74
$ws.Range("H76").Value = 3270846.8
$ws.Range("J76").Value = 6175539.5
$ws.Range("L76").Value = 6175539.5
$ws.Range("N76").Value = -6176169.5
$ws.Range("H79").Value = 3270846.8
$ws.Range("J79").Value = 6175539.5
$ws.Range("L79").Value = 6175539.5
$ws.Range("N79").Value = -6177723.5
$ws.Range("H106").Value = 3169.1724
$ws.Range("I106").Value = 2336.5
$ws.Range("K106").Value = 2336.5
$ws.Range("M106").Value = -1705.5
$ws.Range("H129").Value = 271198.62
$ws.Range("J129").Value = 295104.47
$ws.Range("L129").Value = 885313.4099999999
$ws.Range("N129").Value = -895313.4099999999
$ws.Range("H137").Value = 88566.44
$ws.Range("I137").Value = 122520.336
$ws.Range("J137").Value = 2375.7693
$ws.Range("K137").Value = 367561.008
$ws.Range("L137").Value = 7127.3079
$ws.Range("M137").Value = -365011.008
$ws.Range("N137").Value = -12227.3079
$ws.Range("H138").Value = 5197.4375
$ws.Range("I138").Value = 5179.8
$ws.Range("J138").Value = 5205.4546
$ws.Range("K138").Value = 15539.4
$ws.Range("L138").Value = 15616.3638
$ws.Range("M138").Value = -10399.4
$ws.Range("N138").Value = -25896.3638
$ws.Range("H141").Value = 1841.2632
$ws.Range("I141").Value = 1710.3611
$ws.Range("J141").Value = 4197.5
$ws.Range("K141").Value = 5131.0833
$ws.Range("L141").Value = 12592.5
$ws.Range("M141").Value = 48.91669999999976
$ws.Range("N141").Value = -22952.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1522.375
$ws.Range("I110").Value = 1096.0714
$ws.Range("J110").Value = 4506.5
$ws.Range("K110").Value = 1096.0714
$ws.Range("L110").Value = 4506.5
$ws.Range("M110").Value = 948.9286
$ws.Range("N110").Value = -8596.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 5008
$ws.Range("I24").Value = 1016
$ws.Range("J24").Value = 9000
$ws.Range("K24").Value = 1016
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = -781
$ws.Range("N24").Value = -9470
$ws.Range("H36").Value = 677.6667
$ws.Range("I36").Value = 677.6667
$ws.Range("K36").Value = 677.6667
$ws.Range("M36").Value = -143.6667
$ws.Range("H75").Value = 18942.8
$ws.Range("I75").Value = 9904.666999999999
$ws.Range("K75").Value = 9904.666999999999
$ws.Range("M75").Value = -8968.666999999999
$ws.Range("H78").Value = 18942.8
$ws.Range("I78").Value = 9904.666999999999
$ws.Range("K78").Value = 29714.001
$ws.Range("M78").Value = -25034.001
$ws.Range("H99").Value = 1557.1428
$ws.Range("I99").Value = 1633.3334
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 1633.3334
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = -135.3334
$ws.Range("N99").Value = -4096
$ws.Range("H107").Value = 2359.6829
$ws.Range("I107").Value = 2049.6428
$ws.Range("J107").Value = 3027.4614
$ws.Range("K107").Value = 2049.6428
$ws.Range("L107").Value = 3027.4614
$ws.Range("M107").Value = -129.6428000000001
$ws.Range("N107").Value = -6867.4614

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1178.3334
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1574
$ws.Range("H31").Value = 6012.3887
$ws.Range("I31").Value = 4818.857
$ws.Range("J31").Value = 6300.483
$ws.Range("K31").Value = 4818.857
$ws.Range("L31").Value = 6300.483
$ws.Range("M31").Value = -4523.857
$ws.Range("N31").Value = -6890.483
$ws.Range("H34").Value = 6012.3887
$ws.Range("I34").Value = 4818.857
$ws.Range("J34").Value = 6300.483
$ws.Range("K34").Value = 4818.857
$ws.Range("L34").Value = 6300.483
$ws.Range("M34").Value = -4616.857
$ws.Range("N34").Value = -6704.483
$ws.Range("H48").Value = 9800
$ws.Range("J48").Value = 9800
$ws.Range("L48").Value = 9800
$ws.Range("N48").Value = -10752
$ws.Range("H113").Value = 1178.3334
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340
$ws.Range("H134").Value = 66667696
$ws.Range("I134").Value = 76924070
$ws.Range("K134").Value = 230772210
$ws.Range("M134").Value = -230769675

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 223.54546
$ws.Range("I18").Value = 162.11111
$ws.Range("K18").Value = 486.33333
$ws.Range("M18").Value = -317.33333
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H122").Value = 1113.2142
$ws.Range("J122").Value = 1190.8462
$ws.Range("L122").Value = 10717.6158
$ws.Range("N122").Value = -15617.6158
$ws.Range("H130").Value = 1690.5454
$ws.Range("I130").Value = 1119.2
$ws.Range("J130").Value = 2166.6667
$ws.Range("K130").Value = 3357.6
$ws.Range("L130").Value = 6500.000100000001
$ws.Range("M130").Value = 1662.4
$ws.Range("N130").Value = -16540.0001
$ws.Range("H131").Value = 720.6799999999999
$ws.Range("J131").Value = 764.98865
$ws.Range("L131").Value = 2294.96595
$ws.Range("N131").Value = -12374.96595
$ws.Range("H137").Value = 16672433
$ws.Range("J137").Value = 19614360
$ws.Range("L137").Value = 58843080
$ws.Range("N137").Value = -58853280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 160800
$ws.Range("I24").Value = 200000
$ws.Range("J24").Value = 4000
$ws.Range("K24").Value = 200000
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = -199827
$ws.Range("N24").Value = -4346
$ws.Range("H51").Value = 30326
$ws.Range("J51").Value = 30326
$ws.Range("L51").Value = 30326
$ws.Range("N51").Value = -31344

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2980.1177
$ws.Range("I7").Value = 2750.923
$ws.Range("J7").Value = 3725
$ws.Range("K7").Value = 2750.923
$ws.Range("L7").Value = 3725
$ws.Range("M7").Value = -2638.923
$ws.Range("N7").Value = -3949
$ws.Range("H22").Value = 5900.5
$ws.Range("H27").Value = 5900.5
$ws.Range("H46").Value = 1133
$ws.Range("H50").Value = 16500
$ws.Range("H54").Value = 18000
$ws.Range("J54").Value = 18000
$ws.Range("L54").Value = 18000
$ws.Range("N54").Value = -19288
$ws.Range("H55").Value = 281.7647
$ws.Range("J55").Value = 250
$ws.Range("L55").Value = 250
$ws.Range("N55").Value = -596
$ws.Range("H68").Value = 2374.25
$ws.Range("I68").Value = 2499.6667
$ws.Range("J68").Value = 1998
$ws.Range("K68").Value = 2499.6667
$ws.Range("L68").Value = 1998
$ws.Range("M68").Value = -1750.6667
$ws.Range("N68").Value = -3496
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H71").Value = 2374.25
$ws.Range("I71").Value = 2499.6667
$ws.Range("J71").Value = 1998
$ws.Range("K71").Value = 12498.3335
$ws.Range("L71").Value = 9990
$ws.Range("M71").Value = -8754.333500000001
$ws.Range("N71").Value = -17478
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H126").Value = 2980.1177
$ws.Range("I126").Value = 2750.923
$ws.Range("J126").Value = 3725
$ws.Range("K126").Value = 8252.769
$ws.Range("L126").Value = 11175
$ws.Range("M126").Value = -5782.769
$ws.Range("N126").Value = -16115
$ws.Range("H132").Value = 2651.889
$ws.Range("I132").Value = 1925.5294
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 5776.5882
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -3246.5882
$ws.Range("N132").Value = -50060
$ws.Range("H140").Value = 47904.25
$ws.Range("J140").Value = 47904.25
$ws.Range("L140").Value = 47904.25
$ws.Range("N140").Value = -58264.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 15000
$ws.Range("J31").Value = 15000
$ws.Range("L31").Value = 15000
$ws.Range("N31").Value = -15696
$ws.Range("H107").Value = 518.8889
$ws.Range("I107").Value = 482.73334
$ws.Range("J107").Value = 699.6667
$ws.Range("K107").Value = 1448.20002
$ws.Range("L107").Value = 2099.0001
$ws.Range("M107").Value = 471.79998
$ws.Range("N107").Value = -5939.0001
$ws.Range("H136").Value = 29414948
$ws.Range("I136").Value = 41668280
$ws.Range("K136").Value = 125004840
$ws.Range("M136").Value = -125002290

